# audiences-process.pptx - "Sync non localizable files" update
#
# 1) Refresh the cached "datetimeFigureOut" footer/date field text
#    (slide master + every slide layout) from 01/02/2023 to 02/08/2023.
# 2) Widen the "Target audience(s) in campaigns" textbox on slide 1 and
#    extend its wording to "... and journeys".

$p = $ppt.ActivePresentation

$oldDate = "01/02/2023"
$newDate = "02/08/2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholder {
    param($container)

    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's Date Placeholder
Update-DatePlaceholder $p.SlideMaster

# Every slide layout's Date Placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L)
}

# Slide 1: widen "TextBox 11" and update its copy
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Target audience(s) in campaigns") {
        $shp.Width = 176.23736
        $shp.TextFrame.TextRange.Text = "Target audience(s) in campaigns and journeys"
    }
}
